# Update cryptocurrency price/volume data as of Tue Dec 12 03:25:42 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "41.847.78"
$ws.Range("E2").Value = "  -1.02%  "

Set-TextValue $ws.Range("D3") "2.238.51"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("E4").Value = "  -0.09%  "

Set-TextValue $ws.Range("D5") "252.38"
$ws.Range("E5").Value = "  +8.72%  "

$ws.Range("E6").Value = "  -0.34%  "

Set-TextValue $ws.Range("D7") "72.07"
$ws.Range("E7").Value = "  +1.49%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +0.13%  "

Set-TextValue $ws.Range("D10") "42.33"
$ws.Range("E10").Value = "  +18.90%  "

Set-TextValue $ws.Range("D11") "0.0978"
$ws.Range("E11").Value = "  -1.61%  "

Set-TextValue $ws.Range("D12") "58.37"
$ws.Range("E12").Value = "  -0.89%  "

$ws.Range("E13").Value = "  +0.91%  "

Set-TextValue $ws.Range("D14") "6.99"
$ws.Range("E14").Value = "  +2.67%  "

Set-TextValue $ws.Range("D15") "2.570.46"
$ws.Range("E15").Value = "  -0.29%  "

Set-TextValue $ws.Range("D16") "15.12"
$ws.Range("E16").Value = "  +1.26%  "

Set-TextValue $ws.Range("D17") "0.863"
$ws.Range("E17").Value = "  -0.68%  "

Set-TextValue $ws.Range("D18") "2.237.89"
$ws.Range("E18").Value = "  -0.14%  "

Set-TextValue $ws.Range("D19") "41.781.52"
$ws.Range("E19").Value = "  -0.78%  "

Set-TextValue $ws.Range("D20") "0.0₃0967"
$ws.Range("E20").Value = "  -1.85%  "

Set-TextValue $ws.Range("D21") "73.48"
$ws.Range("E21").Value = "  -0.61%  "

$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("E23").Value = "  +20.07%  "

Set-TextValue $ws.Range("D24") "235.65"
$ws.Range("E24").Value = "  -0.69%  "

Set-TextValue $ws.Range("D25") "0.999"
$ws.Range("E25").Value = "  -0.06%  "

Set-TextValue $ws.Range("D26") "3.75"
$ws.Range("E26").Value = "  +2.17%  "

Set-TextValue $ws.Range("D27") "2.51"
$ws.Range("E27").Value = "  +6.65%  "

Set-TextValue $ws.Range("D28") "10.30"
$ws.Range("E28").Value = "  +2.72%  "

Set-TextValue $ws.Range("D31") "20.80"
$ws.Range("E31").Value = "  +0.16%  "

$ws.Range("E32").Value = "  +2.01%  "

$ws.Range("E33").Value = "  -1.04%  "

$ws.Range("E34").Value = "  +1.44%  "

Set-TextValue $ws.Range("D35") "0.0727"
$ws.Range("E35").Value = "  +1.12%  "

Set-TextValue $ws.Range("D36") "26.58"
$ws.Range("E36").Value = "  +19.99%  "

Set-TextValue $ws.Range("D37") "4.71"
$ws.Range("E37").Value = "  -2.23%  "

Set-TextValue $ws.Range("D38") "4.10"
$ws.Range("E38").Value = "  +13.83%  "

Set-TextValue $ws.Range("D39") "0.0283"
$ws.Range("E39").Value = "  +6.36%  "

Set-TextValue $ws.Range("D40") "2.30"
$ws.Range("E40").Value = "  +2.40%  "

$ws.Range("E41").Value = "  +0.34%  "

Set-TextValue $ws.Range("D42") "69.50"
$ws.Range("E42").Value = "  +5.19%  "

$ws.Range("E43").Value = "  +12.79%  "

$ws.Range("E44").Value = "  +0.35%  "

Set-TextValue $ws.Range("D45") "11.62"
$ws.Range("E45").Value = "  +14.37%  "

$ws.Range("E46").Value = "  +1.10%  "

Set-TextValue $ws.Range("D47") "4.86"
$ws.Range("E47").Value = "  +8.20%  "

Set-TextValue $ws.Range("D48") "0.102"
$ws.Range("E48").Value = "  +0.91%  "

$ws.Range("E49").Value = "  -0.13%  "

$ws.Range("E50").Value = "  +7.82%  "

$ws.Range("E51").Value = "  +1.43%  "

# Rows 29 and 30 swap order: Toncoin now ranked 29th, Monero 30th
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D29") "2.20"
$ws.Range("E29").Value = "  +2.31%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D30") "171.89"
$ws.Range("E30").Value = "  +2.77%  "
